# Project "Sample Project" resave: on the "Rules" sheet, the rule-name
# cell for the fourth rule row (B11, previously "R40") was changed to "1".
#
# "1" looks numeric, but the original/target cell is a text (shared
# string) cell, not a number. Typing a leading apostrophe is how Excel
# keeps a numeric-looking entry stored as text (quote-prefixed text)
# instead of silently converting it to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("B11").Value = "'1"
